# Applies the "extracted additional Nb-rich data from 10.1016/j.msea.2021.142290"
# edit to Sheet1 of the ULTERA contribute workbook.
#
# Summary of the edit:
#   1) Rows 106-108 (NbMoTaV0.25/0.5/0.75) keep referencing the same DOI string;
#      no value change there (shared-string table just gets re-packed once the
#      old, now-unused "minimum tensile ductility" string disappears).
#   2) Rows 109-113 (existing NbMoTaV rows) get re-labelled from tensile to
#      compressive properties, and their Unit [SI] column gets filled in.
#   3) Rows 114-123 are brand new rows describing a Nb40Ti25Al15V10Ta5Hf3W2
#      BCC+B2 alloy (density, compressive yield stress, UCS and minimum
#      compressive ductility at several temperatures) from the same new DOI.
#   4) The active selection moves to L130 (cosmetic "where was the user last
#      looking" bookkeeping that Excel stores in the sheetView).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 2) Re-label the existing NbMoTaV rows (109-113) from tensile to compressive
#    properties, and populate their Unit [SI] (column L).
# ---------------------------------------------------------------------------

# Row 109: tensile yield stress @298K -> compressive yield stress, unit Pa
$ws.Cells.Item(109, 6).Value = "compressive yield stress"
$ws.Cells.Item(109, 12).Value = "Pa"

# Row 110: tensile yield stress @1273K -> compressive yield stress, unit Pa
$ws.Cells.Item(110, 6).Value = "compressive yield stress"
$ws.Cells.Item(110, 12).Value = "Pa"

# Row 111: tensile ductility @298K -> minimum compressive ductility, unit %
$ws.Cells.Item(111, 6).Value = "minimum compressive ductility"
$ws.Cells.Item(111, 12).Value = "%"

# Row 112: minimum tensile ductility @1273K -> minimum compressive ductility, unit %
$ws.Cells.Item(112, 6).Value = "minimum compressive ductility"
$ws.Cells.Item(112, 12).Value = "%"

# Row 113: UTS @1273K -> UCS, unit Pa
$ws.Cells.Item(113, 6).Value = "UCS"
$ws.Cells.Item(113, 12).Value = "Pa"

# ---------------------------------------------------------------------------
# 3) New rows 114-123: Nb40Ti25Al15V10Ta5Hf3W2 (BCC+B2, B2 nanoprecipitates,
#    AAM processing, EXP source, DOI 10.1016/j.msea.2021.142290)
# ---------------------------------------------------------------------------

$composition = "Nb40Ti25Al15V10Ta5Hf3W2"
$structure = "BCC+B2"
$processing = "AAM"
$comment = "B2 nanoprecipitates"
$source = "EXP"
$strainRate = "strain rate 1e-3/s"
$doi = "10.1016/j.msea.2021.142290"

# Row 114: density @298K = 7340 kg/m^3 (no strain-rate parameter for density)
$ws.Cells.Item(114, 2).Value = $composition
$ws.Cells.Item(114, 3).Value = $structure
$ws.Cells.Item(114, 4).Value = $processing
$ws.Cells.Item(114, 5).Value = $comment
$ws.Cells.Item(114, 6).Value = "density"
$ws.Cells.Item(114, 7).Value = $source
$ws.Cells.Item(114, 9).Value = 298
$ws.Cells.Item(114, 10).Value = 7340
$ws.Cells.Item(114, 12).Value = "kg/m^3"
$ws.Cells.Item(114, 14).Value = $doi

# Row 115: compressive yield stress @298K = 1.024e9 Pa +/- 7e6 Pa
$ws.Cells.Item(115, 2).Value = $composition
$ws.Cells.Item(115, 3).Value = $structure
$ws.Cells.Item(115, 4).Value = $processing
$ws.Cells.Item(115, 5).Value = $comment
$ws.Cells.Item(115, 6).Value = "compressive yield stress"
$ws.Cells.Item(115, 7).Value = $source
$ws.Cells.Item(115, 8).Value = $strainRate
$ws.Cells.Item(115, 9).Value = 298
$ws.Cells.Item(115, 10).Value = 1024000000
$ws.Cells.Item(115, 11).Value = 7000000
$ws.Cells.Item(115, 12).Value = "Pa"
$ws.Cells.Item(115, 14).Value = $doi

# Row 116: compressive yield stress @1073K = 6.11e8 Pa +/- 2.4e7 Pa
$ws.Cells.Item(116, 2).Value = $composition
$ws.Cells.Item(116, 3).Value = $structure
$ws.Cells.Item(116, 4).Value = $processing
$ws.Cells.Item(116, 5).Value = $comment
$ws.Cells.Item(116, 6).Value = "compressive yield stress"
$ws.Cells.Item(116, 7).Value = $source
$ws.Cells.Item(116, 8).Value = $strainRate
$ws.Cells.Item(116, 9).Value = 1073
$ws.Cells.Item(116, 10).Value = 611000000
$ws.Cells.Item(116, 11).Value = 24000000
$ws.Cells.Item(116, 12).Value = "Pa"
$ws.Cells.Item(116, 14).Value = $doi

# Row 117: compressive yield stress @1173K = 4.37e8 Pa +/- 8e6 Pa
$ws.Cells.Item(117, 2).Value = $composition
$ws.Cells.Item(117, 3).Value = $structure
$ws.Cells.Item(117, 4).Value = $processing
$ws.Cells.Item(117, 5).Value = $comment
$ws.Cells.Item(117, 6).Value = "compressive yield stress"
$ws.Cells.Item(117, 7).Value = $source
$ws.Cells.Item(117, 8).Value = $strainRate
$ws.Cells.Item(117, 9).Value = 1173
$ws.Cells.Item(117, 10).Value = 437000000
$ws.Cells.Item(117, 11).Value = 8000000
$ws.Cells.Item(117, 12).Value = "Pa"
$ws.Cells.Item(117, 14).Value = $doi

# Row 118: compressive yield stress @1273K = 2.37e8 Pa +/- 8e6 Pa
$ws.Cells.Item(118, 2).Value = $composition
$ws.Cells.Item(118, 3).Value = $structure
$ws.Cells.Item(118, 4).Value = $processing
$ws.Cells.Item(118, 5).Value = $comment
$ws.Cells.Item(118, 6).Value = "compressive yield stress"
$ws.Cells.Item(118, 7).Value = $source
$ws.Cells.Item(118, 8).Value = $strainRate
$ws.Cells.Item(118, 9).Value = 1273
$ws.Cells.Item(118, 10).Value = 237000000
$ws.Cells.Item(118, 11).Value = 8000000
$ws.Cells.Item(118, 12).Value = "Pa"
$ws.Cells.Item(118, 14).Value = $doi

# Row 119: UCS @1073K = 7.46e8 Pa
$ws.Cells.Item(119, 2).Value = $composition
$ws.Cells.Item(119, 3).Value = $structure
$ws.Cells.Item(119, 4).Value = $processing
$ws.Cells.Item(119, 5).Value = $comment
$ws.Cells.Item(119, 6).Value = "UCS"
$ws.Cells.Item(119, 7).Value = $source
$ws.Cells.Item(119, 8).Value = $strainRate
$ws.Cells.Item(119, 9).Value = 1073
$ws.Cells.Item(119, 10).Value = 746000000
$ws.Cells.Item(119, 12).Value = "Pa"
$ws.Cells.Item(119, 14).Value = $doi

# Row 120: UCS @1173K = 4.43e8 Pa
$ws.Cells.Item(120, 2).Value = $composition
$ws.Cells.Item(120, 3).Value = $structure
$ws.Cells.Item(120, 4).Value = $processing
$ws.Cells.Item(120, 5).Value = $comment
$ws.Cells.Item(120, 6).Value = "UCS"
$ws.Cells.Item(120, 7).Value = $source
$ws.Cells.Item(120, 8).Value = $strainRate
$ws.Cells.Item(120, 9).Value = 1173
$ws.Cells.Item(120, 10).Value = 443000000
$ws.Cells.Item(120, 12).Value = "Pa"
$ws.Cells.Item(120, 14).Value = $doi

# Row 121: UCS @1273K = 2.44e8 Pa
$ws.Cells.Item(121, 2).Value = $composition
$ws.Cells.Item(121, 3).Value = $structure
$ws.Cells.Item(121, 4).Value = $processing
$ws.Cells.Item(121, 5).Value = $comment
$ws.Cells.Item(121, 6).Value = "UCS"
$ws.Cells.Item(121, 7).Value = $source
$ws.Cells.Item(121, 8).Value = $strainRate
$ws.Cells.Item(121, 9).Value = 1273
$ws.Cells.Item(121, 10).Value = 244000000
$ws.Cells.Item(121, 12).Value = "Pa"
$ws.Cells.Item(121, 14).Value = $doi

# Row 122: minimum compressive ductility @298K = 70%
$ws.Cells.Item(122, 2).Value = $composition
$ws.Cells.Item(122, 3).Value = $structure
$ws.Cells.Item(122, 4).Value = $processing
$ws.Cells.Item(122, 5).Value = $comment
$ws.Cells.Item(122, 6).Value = "minimum compressive ductility"
$ws.Cells.Item(122, 7).Value = $source
$ws.Cells.Item(122, 8).Value = $strainRate
$ws.Cells.Item(122, 9).Value = 298
$ws.Cells.Item(122, 10).Value = 70
$ws.Cells.Item(122, 12).Value = "%"
$ws.Cells.Item(122, 14).Value = $doi

# Row 123: minimum compressive ductility @1273K = 70%
$ws.Cells.Item(123, 2).Value = $composition
$ws.Cells.Item(123, 3).Value = $structure
$ws.Cells.Item(123, 4).Value = $processing
$ws.Cells.Item(123, 5).Value = $comment
$ws.Cells.Item(123, 6).Value = "minimum compressive ductility"
$ws.Cells.Item(123, 7).Value = $source
$ws.Cells.Item(123, 8).Value = $strainRate
$ws.Cells.Item(123, 9).Value = 1273
$ws.Cells.Item(123, 10).Value = 70
$ws.Cells.Item(123, 12).Value = "%"
$ws.Cells.Item(123, 14).Value = $doi

# ---------------------------------------------------------------------------
# 4) Move the active selection to L130, matching the author's last cursor
#    position when they saved the workbook.
# ---------------------------------------------------------------------------
$ws.Range("L130").Select()
